$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the test steps text in C3 to prepend a new instruction line
$ws.Range("C3").Value = "1) Use web resource -  https://devexpress.github.io/testcafe/example/`n2) Insert Name `n2) Click submit button`n3) Verify Name in message"

# Widen column C to fit the longer text
$ws.Columns("C").ColumnWidth = 58.7109375

# Update the last active selection to C11
$ws.Range("C11").Select()
